$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ("STT" = 2): fill in the "Công việc đã đạt được" (D3) description
# that was previously left blank, now that this week's work (reset-password,
# my-ticket UI, profile page) is finished.
$newWork = @"
1. Hoàn thành mô hình dữ liệu mức ý niệm (ERD) thể hiện các thực thể và mối quan hệ chính của hệ thống bán vé sự kiện, bao gồm Người dùng, Sự kiện, Vé, Đơn hàng và Phân quyền người dùng.
2. Hoàn thiện Use Case Đăng ký, Đăng nhập và Quản lý thông tin cá nhân, mô tả đầy đủ các bước tương tác giữa người dùng và hệ thống. Đảm bảo các Use Case thể hiện rõ phạm vi chức năng, luồng hoạt động và mối liên kết giữa người dùng với hệ thống.
3. Thiết kế và xây dựng giao diện người dùng cho các màn hình đăng ký, đăng nhập và quản lý thông tin cá nhân, đảm bảo bố cục hợp lý, dễ sử dụng và đồng nhất phong cách. Hoàn thiện các chức năng cơ bản cho phép người dùng đăng ký tài khoản, đăng nhập hệ thống, chỉnh sửa và lưu thông tin cá nhân. Thực hiện phân quyền người dùng, đảm bảo mỗi nhóm đối tượng chỉ truy cập được vào các chức năng phù hợp.
Kết quả đạt được:
- Hoàn thành toàn bộ Use Case và giao diện liên quan đến đăng ký, đăng nhập và quản lý thông tin cá nhân, đáp ứng đúng yêu cầu trong phạm vi công việc tuần này, là nền tảng để phát triển các chức năng nâng cao trong những tuần tiếp theo như tạo sự kiện, mua vé và quản lý đơn hàng.
"@

# Match the formatting used by the row above (C2:D2 - left/top aligned,
# vertical-centered, wrapped text) for the now-filled C3:D3 cells.
$ws.Range("C2:D2").Copy()
$ws.Range("C3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D3").Value = $newWork

# Row grows tall enough to show the whole paragraph.
$ws.Rows.Item(3).RowHeight = 303

# --- View state: scrolled down one row, with C2 as the active selected cell.
$ws.Range("C2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
